# Apply updates described by the diff:
# - Rename sheet "Through 2022-03-12" -> "Through 2022-03-13"
# - Update header text "2022 (through 03-12)" -> "2022 (through 03-13)"
# - Update March value in I4: 56 -> 62
# - Update Total value in I14: 356 -> 362

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (this updates the sheet's name in workbook.xml)
$ws.Name = "Through 2022-03-13"

# Update the header cell text (shared string) in column I, row 1
$ws.Range("I1").Value = "2022 (through 03-13)"

# Update the March data value
$ws.Range("I4").Value = 62

# Update the Total row value
$ws.Range("I14").Value = 362
